# Delete the "Peach et al., 2018 / NanoBret" data row from the
# VEGFA165_NRP1 (VEGF:NRP1) sheet.
#
# In the original sheet, row 7 holds:
#   A7 = "Peach et al., 2018"   B7 = "NanoBret"   C7 = 4.95   D7 = 1.25
# and row 8 holds the "Unpublished data" / SPR summary row (with
# AVERAGE/STDEVA formulas referencing G2:G3). Deleting row 7 entirely
# shifts row 8 up to become the new row 7, shrinking the sheet from
# A1:G8 to A1:G7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the whole row so everything below shifts up (formulas, shared
# strings and the sheet dimension all adjust automatically).
$ws.Rows.Item(7).Delete()

# Match the author's resulting selection: the (now last) summary row.
$ws.Range("A7:D7").Select()
